$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New control-signal column: "ctrl_alu_dmem" (header, bold like the rest
# of row 1)
$ws.Range("O1").Value = "ctrl_alu_dmem"
$ws.Range("O1").Font.Bold = $true

# Fill in the ctrl_alu_dmem values for every instruction row (2-17).
# Only "lw" (row 9) asserts this new control signal.
$values = @(0,0,0,0,0,0,0,1,0,0,0,0,0,0,0,0)
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 15).Value = $values[$i]
}

# Give column O an explicit width (best-fit sized to the header text).
$ws.Columns.Item(15).ColumnWidth = 13.666666666666666

# Move the active selection onto the newly added column, which also
# clears the sheet's old frozen/scrolled "topLeftCell" view state.
$ws.Range("O11").Select()
